# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data block (rows 8-9),
# pushing the existing rows 8-27 down to rows 10-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 8 (shifts rows 8:27 -> 10:29)
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# --- New row 8 ---
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44973
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100103
$ws.Range("H8").Value = "Frutos de hueso (carozo)"
$ws.Range("I8").Value = 100103002
$ws.Range("J8").Value = "Ciruela"
$ws.Range("K8").Value = "Larry Ann"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19000
$ws.Range("Q8").Value = "$/bandeja 18 kilos granel"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1056
$ws.Range("T8").Value = 18

# --- New row 9 ---
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44973
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100103
$ws.Range("H9").Value = "Frutos de hueso (carozo)"
$ws.Range("I9").Value = 100103002
$ws.Range("J9").Value = "Ciruela"
$ws.Range("K9").Value = "Pink Delight"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 270
$ws.Range("N9").Value = 19000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 19500
$ws.Range("Q9").Value = "$/bandeja 18 kilos granel"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1083
$ws.Range("T9").Value = 18
